# Add a "GEN_ROW" flag column to Sheet1 so the generator knows which rows to
# emit. This inserts a new column B (pushing the existing columns B:L to
# C:M), then fills in the header and the Y/N flag values for the sample
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank column at B; Excel copies formatting from the column
# to its left (A), which is exactly the style pattern the target file uses
# for the new column.
$ws.Columns.Item(2).EntireColumn.Insert()

# New column width matches column A's width (16.5 chars); header/values follow.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()

$ws.Range("B1").Value = "GEN_ROW"
$ws.Range("B3").Value = "Y"
$ws.Range("B2").Value = "N"
# Row 4 is intentionally left blank in column B.

# Re-apply the number formats/styles from column A onto the new column B
# data cells so they keep the quote-prefixed "text" style used for the
# UNIQUE_ID_COLUMN values instead of the default.
$ws.Range("A2:A4").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the author's recorded selection after the edit.
[void]$ws.Range("B4").Select()
